# Simulate re-running the portfolio Monte Carlo / DCF simulation (dcf_sim)
# which recomputes the realized "Capital out" proceeds for each company
# and refreshes the resulting MOIC / XIRR figures downstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")

# New simulated capital-out results for companies 1 and 2 (column D).
$ws.Range("D16").Value = 5
$ws.Range("D17").Value = 5

# Leave selection where the user last clicked after reviewing the results.
$ws.Range("E16").Select()
